# Rename the "_old"/"_new" column header suffixes to the respective
# format-version suffixes ("_FV2410" for the old/left block, "_FV2504"
# for the new/right block), then turn the sheet's data range into a
# native Excel table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2504"
}

# Turn A1:U90 into a proper Excel table (adds autofilter + tableParts).
$range = $ws.Range("A1:U90")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1, i.e. top-left cell of the
# scrolling pane is A2).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
